$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "2025-03-27"
$ws.Range("B34").Value = "développement durable"
$ws.Range("C34").Value = 60
$ws.Range("D34").Value = 1
$ws.Range("A34:D34").ClearFormats()
